$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SimCLR")
$ws.Range("B3").Copy()
$ws.Range("A16:C20").PasteSpecial(-4122)
$ws.Range("A16:C20").Copy()
$ws.Range("A25:C29").PasteSpecial(-4122)
